$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '29.878.74'
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = '1.888.60'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '0.7703'
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("D6").Value = '242.73'
$ws.Range("E6").Value = '  -0.88%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '0.3120'
$ws.Range("E8").Value = '  -0.78%  '
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").Value = '0.07179'
$ws.Range("E10").Value = '  -1.62%  '
$ws.Range("D11").Value = '0.08615'
$ws.Range("E11").Value = '  +6.29%  '
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").Value = '0.7637'
$ws.Range("E12").Value = '  -0.91%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.941.80'
$ws.Range("E13").Value = '  +2.63%  '
$ws.Range("D14").Value = '5.368'
$ws.Range("E14").Value = '  -2.75%  '
$ws.Range("D15").Value = '93.62'
$ws.Range("E15").Value = '  +1.27%  '
$ws.Range("D16").Value = '6.183'
$ws.Range("E16").Value = '  -2.53%  '
$ws.Range("D17").Value = '29.975.91'
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("D18").Value = '13.78'
$ws.Range("E18").Value = '  -1.41%  '
$ws.Range("D19").Value = '244.52'
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").Value = '0.000007815'
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").Value = '2.235.82'
$ws.Range("E21").Value = '  +4.52%  '
$ws.Range("D22").Value = '1.0000'
$ws.Range("D23").Value = '8.054'
$ws.Range("E23").Value = '  -1.30%  '
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("D25").Value = '0.1660'
$ws.Range("E25").Value = '  +3.68%  '
$ws.Range("D26").Value = '9.368'
$ws.Range("E26").Value = '  -1.21%  '
$ws.Range("D27").Value = '162.26'
$ws.Range("E27").Value = '  -0.27%  '
$ws.Range("E28").Value = '  -0.17%  '
$ws.Range("D29").Value = '2.040'
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("D30").Value = '1.454'
$ws.Range("E30").Value = '  +0.90%  '
$ws.Range("D31").Value = '1.534'
$ws.Range("E31").Value = '  -1.11%  '
$ws.Range("D32").Value = '4.504'
$ws.Range("E32").Value = '  +0.39%  '
$ws.Range("D33").Value = '4.103'
$ws.Range("E33").Value = '  +0.26%  '
$ws.Range("D34").Value = '0.05448'
$ws.Range("E34").Value = '  -1.76%  '
$ws.Range("D35").Value = '1.240'
$ws.Range("E35").Value = '  -1.62%  '
$ws.Range("D36").Value = '0.7450'
$ws.Range("E36").Value = '  -1.33%  '
$ws.Range("D37").Value = '1.003'
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("E38").Value = '  +2.35%  '
$ws.Range("E39").Value = '  +1.73%  '
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("D41").Value = '0.4466'
$ws.Range("E41").Value = '  +0.63%  '
$ws.Range("D42").Value = '1.107.99'
$ws.Range("E42").Value = '  -5.27%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '73.16'
$ws.Range("E43").Value = '  -1.25%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '6.081'
$ws.Range("E44").Value = '  +2.45%  '
$ws.Range("D45").Value = '0.8510'
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("E46").Value = '  +0.11%  '
$ws.Range("D47").Value = '102.46'
$ws.Range("E47").Value = '  -0.15%  '
$ws.Range("D48").Value = '1.867'
$ws.Range("E48").Value = '  -1.58%  '
$ws.Range("D49").Value = '7.658'
$ws.Range("E49").Value = '  +2.33%  '
$ws.Range("D50").Value = '2.121.30'
$ws.Range("E50").Value = '  +3.42%  '
$ws.Range("D51").Value = '2.982'
$ws.Range("E51").Value = '  -2.78%  '
